$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.122.56"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.277.50"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.16"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.83"
$ws.Range("E6").Value = "  -7.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "3.272.63"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -4.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.569"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.77"
$ws.Range("E12").Value = "  -4.55%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "685.73"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").Value = "3.798.27"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.22"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").Value = "67.133.18"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "3.265.57"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.63"
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.879"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.82"
$ws.Range("E23").Value = "  -5.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.19"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.78"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.81"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.01"
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.58"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "581.65"
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.78"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "3.814.01"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.46"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  -16.26%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.38"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.56"
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.31"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("D43").Value = "0.0₃0654"
$ws.Range("E43").Value = "  -6.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.96"
$ws.Range("E44").Value = "  -7.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.323"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0400"
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.52"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.34"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.27"
$ws.Range("E51").Value = "  -0.60%  "
